$d = $word.ActiveDocument

$pairs = @(
    @("28×75=2100", "26×84=2184"),
    @("80×39=3120", "21×26=546"),
    @("84×70=5880", "91×17=1547"),
    @("50×78=3900", "19×74=1406"),
    @("98×32=3136", "42×12=504"),
    @("39×91=3549", "47×95=4465"),
    @("27×32=864", "11×92=1012"),
    @("74×89=6586", "80×62=4960"),
    @("63×14=882", "92×32=2944"),
    @("20×64=1280", "30×20=600"),
    @("69×54=3726", "82×51=4182"),
    @("82×20=1640", "20×34=680"),
    @("98×94=9212", "79×80=6320"),
    @("57×82=4674", "12×41=492"),
    @("36×46=1656", "91×61=5551"),
    @("35×64=2240", "65×85=5525"),
    @("78×23=1794", "43×29=1247"),
    @("49×92=4508", "30×86=2580"),
    @("64×26=1664", "61×58=3538"),
    @("29×99=2871", "75×58=4350"),
    @("98×36=3528", "83×73=6059"),
    @("16×83=1328", "94×20=1880"),
    @("11×53=583", "39×69=2691"),
    @("47×92=4324", "70×52=3640"),
    @("95×55=5225", "14×89=1246")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
